# Update A2: "Nueva ver" -> "Test1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Test1"

# Remove row 3 ("Nueva ver v1000" / "v1000") entirely - shrinks the used range to A1:B2
$ws.Rows(3).Delete()

# Move the active selection to B12, matching the saved cursor position
$ws.Range("B12").Select()
